# Fix list level numbering: paragraphs that were incorrectly bumped to the
# second outline level (lvl="1", i.e. COM IndentLevel 2) when a list style
# was applied should actually sit at the top level (lvl="0", IndentLevel 1).
# This walks every slide/shape/paragraph in the deck and pulls back any
# paragraph sitting at IndentLevel 2 to IndentLevel 1, leaving every other
# paragraph (already at level 1, or explicit continuations deeper than
# level 2) untouched.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)

        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count

            for ($pi = 1; $pi -le $paraCount; $pi++) {
                $para = $tr.Paragraphs($pi, 1)

                if ($para.IndentLevel -eq 2) {
                    $para.IndentLevel = 1
                }
            }
        }
    }
}
